$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.288.65"
$ws.Range("E2").Value = "  +2.04%  "
$ws.Range("D3").Value = "1.895.33"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.41%  "
$ws.Range("E6").Value = "  -0.14%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5179"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4024"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08414"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.36%  "
$ws.Range("E11").Value = "  -0.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.22"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.441"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.07%  "
$ws.Range("D14").Value = "1.898.51"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.319"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "94.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06648"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.52%  "
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("E21").Value = "  -0.12%  "
$ws.Range("E22").Value = "  -1.51%  "
$ws.Range("D23").Value = "30.273.74"
$ws.Range("E23").Value = "  +1.98%  "
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.228"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.64%  "
$ws.Range("D26").Value = "2.106.53"
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.65"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.68"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("E29").Value = "  -4.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "129.47"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.92%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.091"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.23%  "
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.102"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.740"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.29%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02493"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06553"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.300"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.55%  "
$ws.Range("E38").Value = "  -0.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.219"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.05%  "
$ws.Range("E40").Value = "  +4.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.774"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6499"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.02%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.227"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6092"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.78%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.687"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.054"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("E48").Value = "  -0.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.53"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.13%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.157"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.13"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.16%  "
